# Weekly update: a new price record (week of 2021-11-10) is inserted at the
# top of the data table (row 30), pushing the existing rows 30-47 down to
# 31-48. Using a real row insert (instead of rewriting every shifted row by
# hand) lets Excel carry the existing formatting (e.g. the date style on
# column D) down automatically, matching the original author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 30 - everything from the old
# row 30 downward (through row 47) shifts down by one, to rows 31-48.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new week's record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44510
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100112026
$ws.Range("G30").Value = "Haba"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 8000
$ws.Range("N30").Value = '$/saco 25 kilos'
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 320
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
